$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- New header cells (write E1 before D1 so shared-string order matches
#     the target: ADC3 -> index 2, Vcc(-) -> index 3) ---
$ws.Range("E1").Value = "ADC3"
$ws.Range("D1").Value = "Vcc(-)"

# --- Right-align the header row cells (new cellXf, applyAlignment) ---
$ws.Range("A1:B1").HorizontalAlignment = -4152
$ws.Range("D1:E1").HorizontalAlignment = -4152

# --- New ADC3 data (columns D/E), written before the new chart is created
#     so ChartObjects().Add / SetSourceData can see real values ---
$ws.Range("D2").Value = -16
$ws.Range("D3").Value = -15
$ws.Range("D4").Value = -14
$ws.Range("D5").Value = -13
$ws.Range("D6").Value = -12
$ws.Range("D7").Value = -11

$ws.Range("E2").Value = 1018
$ws.Range("E3").Value = 953
$ws.Range("E4").Value = 891
$ws.Range("E5").Value = 833
$ws.Range("E6").Value = 767
$ws.Range("E7").Value = 706

# --- New summary row with formulas ---
$ws.Range("A10").Formula = "=(62.686+62.171)/2"
$ws.Range("B10").Formula = "=(18.41+22.019)/2"

# --- Selection moves to B11 ---
$ws.Range("B11").Select() | Out-Null

# --- Add the new ADC3 chart (Chart 2) below the first one, *before* touching
#     the existing chart, so the engine doesn't spawn a stray extra part ---
$chart2Obj = $ws.ChartObjects().Add(350.625, 225.75, 433.0625, 216)
$chart2Obj.Name = "Chart 2"
$chart2 = $chart2Obj.Chart
$chart2.ChartType = -4169
$chart2.SetSourceData($ws.Range("E1:E7"))
$s2 = $chart2.SeriesCollection(1)
$s2.XValues = $ws.Range("D2:D7")
$s2.Smooth = $true
$t2 = $s2.Trendlines().Add()
$t2.Type = -4132
$t2.DisplayEquation = $true
$chart2.HasTitle = $false
$chart2.HasLegend = $true
$chart2.Legend.Position = -4152

# --- Now move/resize the existing chart (Chart 3) and refresh its cached
#     ADC2 values ---
$chart1Obj = $ws.ChartObjects().Item(1)
$chart1Obj.Top = 0.75
$chart1Obj.Left = 350.625
$chart1Obj.Width = 433.0625
$chart1Obj.Height = 216

# --- Updated ADC2 (column B) readings ---
$ws.Range("B2").Value = 1022
$ws.Range("B3").Value = 959
$ws.Range("B4").Value = 894
$ws.Range("B5").Value = 834
$ws.Range("B6").Value = 771
$ws.Range("B7").Value = 708

$chart1 = $chart1Obj.Chart
$s1 = $chart1.SeriesCollection(1)
$s1.Values = $ws.Range("B2:B7")
